# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Damasco" (Castle Brite, 18 kilos granel)
# right above the former row 21, pushing the existing rows 21-35 down to 24-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 21 (each Insert() pushes rows 21.. down by one)
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# Common, constant columns shared by every Damasco / Vega Central Mapocho row
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103003
$categoria   = "Damasco"
$origen      = "Provincia de San Felipe de Aconcagua"

# --- Row 21 ---
$r = 21
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44539
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Castle Brite"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 290
$ws.Cells.Item($r, 14).Value = 18000
$ws.Cells.Item($r, 15).Value = 18000
$ws.Cells.Item($r, 16).Value = 18000
$ws.Cells.Item($r, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1000
$ws.Cells.Item($r, 20).Value = 18

# --- Row 22 ---
$r = 22
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44539
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Castle Brite"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 16000
$ws.Cells.Item($r, 15).Value = 16000
$ws.Cells.Item($r, 16).Value = 16000
$ws.Cells.Item($r, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 889
$ws.Cells.Item($r, 20).Value = 18

# --- Row 23 ---
$r = 23
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44539
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Castle Brite"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 300
$ws.Cells.Item($r, 14).Value = 14000
$ws.Cells.Item($r, 15).Value = 14000
$ws.Cells.Item($r, 16).Value = 14000
$ws.Cells.Item($r, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 778
$ws.Cells.Item($r, 20).Value = 18

# Make sure the date cells keep the existing date number format used by the
# rest of column D (re-applying Value can clear formatting picked up from
# the row-insert, so set it explicitly from a known-good neighbour cell).
$dateFormat = $ws.Cells.Item(24, 4).NumberFormat
$ws.Cells.Item(21, 4).NumberFormat = $dateFormat
$ws.Cells.Item(22, 4).NumberFormat = $dateFormat
$ws.Cells.Item(23, 4).NumberFormat = $dateFormat
